$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: replace the "Integer" example with a "List" example -----------
$ws.Cells.Item(4, 1).Value  = "hasList"     # A4  name
$ws.Cells.Item(4, 2).Value  = "hasList"     # B4  label_en
$ws.Cells.Item(4, 3).Value  = "List"        # C4  label_de
$ws.Cells.Item(4, 7).Value  = "list"        # G4  comment_en
$ws.Cells.Item(4, 8).Value  = "list"        # H4  comment_de
$ws.Cells.Item(4, 13).Value = "ListValue"   # M4  object
$ws.Cells.Item(4, 14).Value = "List"        # N4  gui_element
# L4 (super = hasValue) and O4 (gui_attributes) keep their existing values.

# --- Drop the now-superfluous blank placeholder cells ----------------------
# (empty cells that only carried a no-op "general alignment" style)
$ws.Range("D2:F2").Clear()
$ws.Range("I2").Clear()
$ws.Range("K2").Clear()

$ws.Range("D3:F3").Clear()
$ws.Range("I3").Clear()
$ws.Range("K3").Clear()

$ws.Range("D4:F4").Clear()
$ws.Range("I4").Clear()
$ws.Range("K4").Clear()

$ws.Range("A5").Clear()
$ws.Range("C5:I5").Clear()
$ws.Range("K5:O5").Clear()
